# "field work and constitution final"
#
# 1) Update the cached datetimeFigureOut field text ("1/23/2019" ->
#    "1/31/2019") on the slide master and on every slide layout's
#    "Date Placeholder" shape.
# 2) Delete the "Secretary" textbox ("TextBox 6") and the straight
#    connector that ran into it ("Straight Connector 24") from slide 1.

$p = $ppt.ActivePresentation

$oldDate = "1/23/2019"
$newDate = "1/31/2019"

# --- 1) Refresh the cached date field text everywhere it appears ----------

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout (CustomLayout) hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- 2) Remove the "Secretary" box and its connector on slide 1 -----------

$s = $p.Slides.Item(1)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Straight Connector 24") {
        $sh.Delete()
    }
}

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 6") {
        $sh.Delete()
    }
}
